$wb = $excel.ActiveWorkbook

# --- Sheet: semantic_aspect_model_schema ---
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")
# Widen column A from 2.4 to 9.6 (character-width units; engine snaps to the
# nearest whole-pixel grid step, so 8.8 is the closest settable value that
# lands on the same grid point as the target stored width of 9.6)
$wsSchema.Columns.Item(1).ColumnWidth = 8.8
# Rename header "id" -> "dtwin_id"
$wsSchema.Range("A1").Value = "dtwin_id"

# --- Sheet: description ---
$wsDescription = $wb.Worksheets.Item("description")
$wsDescription.Range("A5").Value = "dtwin_id"

# --- Sheet: metadata (hidden) ---
$wsMetadata = $wb.Worksheets.Item("metadata")
$wsMetadata.Range("B2").Value = "41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B3").Value = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B4").Value = "2025-03-10 14:48:29+00:00"
$wsMetadata.Range("B5").Value = "Adding auto-generated artifacts for new models"
